$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''62.369.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.77%  '

$ws.Range("D3").Value = '''2.424.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").Value = '''563.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.08%  '

$ws.Range("D6").Value = '''142.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.20%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("D8").Value = '''0.529'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.31%  '

$ws.Range("D9").Value = '''2.419.93'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.68%  '

$ws.Range("D10").Value = '''0.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.18%  '

$ws.Range("E11").Value = '  +1.12%  '

$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").Value = '''5.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.61%  '

$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").Value = '''0.350'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.60%  '

$ws.Range("D14").Value = '''26.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.98%  '

$ws.Range("D15").Value = '''0.0000172'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.81%  '

$ws.Range("D16").Value = '''2.862.37'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = '''62.284.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.92%  '

$ws.Range("D18").Value = '''2.422.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.75%  '

$ws.Range("D19").Value = '''11.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.76%  '

$ws.Range("D20").Value = '''7.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.51%  '

$ws.Range("D21").Value = '''323.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("D22").Value = '''4.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.16%  '

$ws.Range("D23").Value = '''2.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.11%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("D25").Value = '''64.83'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.25%  '

$ws.Range("D26").Value = '''612.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.08%  '

$ws.Range("D27").Value = '''8.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.73%  '

$ws.Range("D28").Value = '''0.0₃0961'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.19%  '

$ws.Range("E29").Value = '  -3.87%  '

$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.36%  '

$ws.Range("D31").Value = '''1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.39%  '

$ws.Range("D32").Value = '''8.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.74%  '

$ws.Range("E33").Value = '  -2.96%  '

$ws.Range("D34").Value = '''0.134'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.16%  '

$ws.Range("D35").Value = '''4.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.70%  '

$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Value = '''1.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.53%  '

$ws.Range("D38").Value = '''0.373'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.46%  '

$ws.Range("D39").Value = '''18.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.20%  '

$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '''5.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.14%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''146.35'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.96%  '

$ws.Range("D42").Value = '''1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.37%  '

$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '''42.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.88%  '

$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.77%  '

$ws.Range("D46").Value = '''143.84'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.63%  '

$ws.Range("D47").Value = '''3.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.73%  '

$ws.Range("D48").Value = '''20.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.32%  '

$ws.Range("D49").Value = '''0.0522'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.93%  '

$ws.Range("D50").Value = '''0.591'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.02%  '

$ws.Range("D51").Value = '''0.0227'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.50%  '
